$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.204.99"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").Value = "'2.912.49"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'369.43"
$ws.Range("E5").Value = "  +5.17%  "
$ws.Range("D6").Value = "'103.35"
$ws.Range("E6").Value = "  -3.21%  "
$ws.Range("E7").Value = "  -2.30%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.591"
$ws.Range("E9").Value = "  -2.79%  "
$ws.Range("D10").Value = "'36.82"
$ws.Range("E10").Value = "  -1.69%  "
$ws.Range("E11").Value = "  +1.64%  "
$ws.Range("E12").Value = "  -1.42%  "
$ws.Range("D13").Value = "'18.47"
$ws.Range("E13").Value = "  -2.25%  "
$ws.Range("D14").Value = "'3.367.00"
$ws.Range("E14").Value = "  -0.44%  "
$ws.Range("D15").Value = "'7.40"
$ws.Range("E15").Value = "  -2.83%  "
$ws.Range("D16").Value = "'2.910.73"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("D17").Value = "'0.944"
$ws.Range("E17").Value = "  -1.14%  "
$ws.Range("D18").Value = "'51.138.83"
$ws.Range("E18").Value = "  -1.07%  "
$ws.Range("D19").Value = "'3.26"
$ws.Range("E19").Value = "  -3.94%  "
$ws.Range("E20").Value = "  -1.68%  "
$ws.Range("D21").Value = "'12.80"
$ws.Range("E21").Value = "  -3.97%  "
$ws.Range("E22").Value = "  -1.08%  "
$ws.Range("D23").Value = "'68.42"
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("D24").Value = "'260.38"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("E25").Value = "  -1.37%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").Value = "'25.75"
$ws.Range("E28").Value = "  -2.23%  "
$ws.Range("D29").Value = "'7.02"
$ws.Range("E29").Value = "  -5.06%  "
$ws.Range("E30").Value = "  -0.53%  "
$ws.Range("D31").Value = "'9.93"
$ws.Range("E31").Value = "  -2.35%  "
$ws.Range("D32").Value = "'6.06"
$ws.Range("E32").Value = "  +2.81%  "
$ws.Range("E33").Value = "  -1.33%  "
$ws.Range("D34").Value = "'34.73"
$ws.Range("E34").Value = "  -2.18%  "
$ws.Range("D35").Value = "'50.91"
$ws.Range("E35").Value = "  -1.17%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.40%  "
$ws.Range("D37").Value = "'0.0421"
$ws.Range("E37").Value = "  -0.52%  "
$ws.Range("D38").Value = "'3.04"
$ws.Range("E38").Value = "  -2.36%  "
$ws.Range("D39").Value = "'2.66"
$ws.Range("E39").Value = "  +1.13%  "
$ws.Range("D40").Value = "'17.11"
$ws.Range("E40").Value = "  -2.39%  "
$ws.Range("E41").Value = "  -5.17%  "
$ws.Range("E42").Value = "  -1.88%  "
$ws.Range("D43").Value = "'22.08"
$ws.Range("E43").Value = "  -1.86%  "
$ws.Range("D44").Value = "'119.69"
$ws.Range("E44").Value = "  +1.24%  "
$ws.Range("D45").Value = "'2.07"
$ws.Range("E45").Value = "  -3.89%  "
$ws.Range("D46").Value = "'2.019.01"
$ws.Range("E46").Value = "  -3.62%  "
$ws.Range("E47").Value = "  -6.26%  "
$ws.Range("E48").Value = "  -3.71%  "
$ws.Range("D49").Value = "'3.186.76"
$ws.Range("E49").Value = "  -0.84%  "
$ws.Range("D50").Value = "'0.238"
$ws.Range("E50").Value = "  +0.96%  "
$ws.Range("E51").Value = "  -7.20%  "
